$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.909.64"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "206.28"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "22.11"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.770.39"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "1.547.86"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "26.911.52"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "217.06"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "9.22"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "154.24"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "14.93"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "1.417.41"
$ws.Range("E33").Value = "  +3.69%  "
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "0.528"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "64.49"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").Value = "1.684.09"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "87.44"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("E51").Value = "  +0.54%  "
